$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.020349605984689
$ws.Cells.Item(2, 4).Value = 1.022174629359064
$ws.Cells.Item(2, 5).Value = 1.029836972115252
$ws.Cells.Item(2, 6).Value = 1.037336229550243
$ws.Cells.Item(2, 9).Value = 1.023594999628091
$ws.Cells.Item(2, 10).Value = 1.02554721649973
$ws.Cells.Item(2, 11).Value = 1.025009786906153
$ws.Cells.Item(2, 12).Value = 1.032649726007894
$ws.Cells.Item(2, 13).Value = 1.040127404791578
$ws.Cells.Item(2, 14).Value = 1.012499121358909
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.021463424500027
$ws.Cells.Item(3, 4).Value = 1.023144942643183
$ws.Cells.Item(3, 5).Value = 1.030864681332056
$ws.Cells.Item(3, 6).Value = 1.038536172828092
$ws.Cells.Item(3, 9).Value = 1.023504579208683
$ws.Cells.Item(3, 10).Value = 1.026297326841434
$ws.Cells.Item(3, 11).Value = 1.025786262692901
$ws.Cells.Item(3, 12).Value = 1.03348507388885
$ws.Cells.Item(3, 13).Value = 1.041136101883477
$ws.Cells.Item(3, 14).Value = 1.012756310177493
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.022184695731987
$ws.Cells.Item(4, 4).Value = 1.023773588947807
$ws.Cells.Item(4, 5).Value = 1.031530533551314
$ws.Cells.Item(4, 6).Value = 1.039313808765389
$ws.Cells.Item(4, 9).Value = 1.023443335729026
$ws.Cells.Item(4, 10).Value = 1.026782745325733
$ws.Cells.Item(4, 11).Value = 1.026288871349019
$ws.Cells.Item(4, 12).Value = 1.034025853867295
$ws.Cells.Item(4, 13).Value = 1.04178942196542
$ws.Cells.Item(4, 14).Value = 1.012922514777585
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.022488051932126
$ws.Cells.Item(5, 4).Value = 1.024038060559989
$ws.Cells.Item(5, 5).Value = 1.031810662680135
$ws.Cells.Item(5, 6).Value = 1.03964101271398
$ws.Cells.Item(5, 9).Value = 1.023416932628374
$ws.Cells.Item(5, 10).Value = 1.026986826549575
$ws.Cells.Item(5, 11).Value = 1.026500210362561
$ws.Cells.Item(5, 12).Value = 1.034253258650684
$ws.Cells.Item(5, 13).Value = 1.042064227717361
$ws.Cells.Item(5, 14).Value = 1.012992335863613
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.022538994635924
$ws.Cells.Item(6, 4).Value = 1.024082477569908
$ws.Cells.Item(6, 5).Value = 1.031857709606911
$ws.Cells.Item(6, 6).Value = 1.039695968466207
$ws.Cells.Item(6, 9).Value = 1.023412460900047
$ws.Cells.Item(6, 10).Value = 1.02702109334905
$ws.Cells.Item(6, 11).Value = 1.026535697595033
$ws.Cells.Item(6, 12).Value = 1.034291444486443
$ws.Cells.Item(6, 13).Value = 1.042110377614219
$ws.Cells.Item(6, 14).Value = 1.013004056116885
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.022188748665937
$ws.Cells.Item(7, 4).Value = 1.023777122089087
$ws.Cells.Item(7, 5).Value = 1.031534275845471
$ws.Cells.Item(7, 6).Value = 1.039318179757125
$ws.Cells.Item(7, 9).Value = 1.023442985510381
$ws.Cells.Item(7, 10).Value = 1.026785472224582
$ws.Cells.Item(7, 11).Value = 1.02629169510493
$ws.Cells.Item(7, 12).Value = 1.034028892221836
$ws.Cells.Item(7, 13).Value = 1.041793093343689
$ws.Cells.Item(7, 14).Value = 1.012923447932227
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.020725911381005
$ws.Cells.Item(8, 4).Value = 1.022502387893734
$ws.Cells.Item(8, 5).Value = 1.030184114099452
$ws.Cells.Item(8, 6).Value = 1.037741509306102
$ws.Cells.Item(8, 9).Value = 1.023565007130014
$ws.Cells.Item(8, 10).Value = 1.025800710116416
$ws.Cells.Item(8, 11).Value = 1.025272163894561
$ws.Cells.Item(8, 12).Value = 1.032931983170499
$ws.Cells.Item(8, 13).Value = 1.040468169287145
$ws.Cells.Item(8, 14).Value = 1.012586083693492
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.018152436874261
$ws.Cells.Item(9, 4).Value = 1.020262189021228
$ws.Cells.Item(9, 5).Value = 1.02781150883515
$ws.Cells.Item(9, 6).Value = 1.034972331664743
$ws.Cells.Item(9, 9).Value = 1.023759139049229
$ws.Cells.Item(9, 10).Value = 1.024065790219161
$ws.Cells.Item(9, 11).Value = 1.02347697498826
$ws.Cells.Item(9, 12).Value = 1.031001031962369
$ws.Cells.Item(9, 13).Value = 1.038138265477068
$ws.Cells.Item(9, 14).Value = 1.011989974284119
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.016439591257137
$ws.Cells.Item(10, 4).Value = 1.018772792231163
$ws.Cells.Item(10, 5).Value = 1.026234173601192
$ws.Cells.Item(10, 6).Value = 1.033132319420464
$ws.Cells.Item(10, 9).Value = 1.023874603872166
$ws.Cells.Item(10, 10).Value = 1.022909411080321
$ws.Cells.Item(10, 11).Value = 1.022281092265943
$ws.Cells.Item(10, 12).Value = 1.029715031074178
$ws.Cells.Item(10, 13).Value = 1.036588195265999
$ws.Cells.Item(10, 14).Value = 1.011591480050356
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.015698565473323
$ws.Cells.Item(11, 4).Value = 1.018128832288674
$ws.Cells.Item(11, 5).Value = 1.025552212411625
$ws.Cells.Item(11, 6).Value = 1.032337016077355
$ws.Cells.Item(11, 9).Value = 1.023921309042974
$ws.Cells.Item(11, 10).Value = 1.022408739773228
$ws.Cells.Item(11, 11).Value = 1.021763476684015
$ws.Cells.Item(11, 12).Value = 1.029158484755441
$ws.Cells.Item(11, 13).Value = 1.035917752253726
$ws.Cells.Item(11, 14).Value = 1.011418671126737
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.015423411878731
$ws.Cells.Item(12, 4).Value = 1.017889780852791
$ws.Cells.Item(12, 5).Value = 1.025299056986841
$ws.Cells.Item(12, 6).Value = 1.032041819855414
$ws.Cells.Item(12, 9).Value = 1.023938164268129
$ws.Cells.Item(12, 10).Value = 1.022222775200687
$ws.Cells.Item(12, 11).Value = 1.021571242525713
$ws.Cells.Item(12, 12).Value = 1.028951803568309
$ws.Cells.Item(12, 13).Value = 1.035668831736869
$ws.Cells.Item(12, 14).Value = 1.0113544434857
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.015482428912444
$ws.Cells.Item(13, 4).Value = 1.017941051683968
$ws.Cells.Item(13, 5).Value = 1.025353352659628
$ws.Cells.Item(13, 6).Value = 1.03210513073432
$ws.Cells.Item(13, 9).Value = 1.023934571058996
$ws.Cells.Item(13, 10).Value = 1.022262664917883
$ws.Cells.Item(13, 11).Value = 1.02161247599071
$ws.Cells.Item(13, 12).Value = 1.028996135357016
$ws.Cells.Item(13, 13).Value = 1.035722220977581
$ws.Cells.Item(13, 14).Value = 1.011368222283054
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.015675819230671
$ws.Cells.Item(14, 4).Value = 1.018109069286104
$ws.Cells.Item(14, 5).Value = 1.02553128333795
$ws.Cells.Item(14, 6).Value = 1.032312610686158
$ws.Cells.Item(14, 9).Value = 1.023922712353276
$ws.Cells.Item(14, 10).Value = 1.022393367741351
$ws.Cells.Item(14, 11).Value = 1.021747585910623
$ws.Cells.Item(14, 12).Value = 1.029141399509654
$ws.Cells.Item(14, 13).Value = 1.035897174126639
$ws.Cells.Item(14, 14).Value = 1.011413362839077
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.015794986131469
$ws.Cells.Item(15, 4).Value = 1.018212609549888
$ws.Cells.Item(15, 5).Value = 1.025640932876699
$ws.Cells.Item(15, 6).Value = 1.032440474383947
$ws.Cells.Item(15, 9).Value = 1.023915340507297
$ws.Cells.Item(15, 10).Value = 1.022473898990777
$ws.Cells.Item(15, 11).Value = 1.02183083573628
$ws.Cells.Item(15, 12).Value = 1.029230907489905
$ws.Cells.Item(15, 13).Value = 1.036004983346332
$ws.Cells.Item(15, 14).Value = 1.011441170298477
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.016488784264058
$ws.Cells.Item(16, 4).Value = 1.01881554992392
$ws.Cells.Item(16, 5).Value = 1.026279454913947
$ws.Cells.Item(16, 6).Value = 1.033185131246975
$ws.Cells.Item(16, 9).Value = 1.023871434968932
$ws.Cells.Item(16, 10).Value = 1.022942639986949
$ws.Cells.Item(16, 11).Value = 1.022315449107855
$ws.Cells.Item(16, 12).Value = 1.029751973530507
$ws.Cells.Item(16, 13).Value = 1.036632706073421
$ws.Cells.Item(16, 14).Value = 1.011602943375499
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.016924158244312
$ws.Cells.Item(17, 4).Value = 1.019194015107871
$ws.Cells.Item(17, 5).Value = 1.026680259955463
$ws.Cells.Item(17, 6).Value = 1.033652618342809
$ws.Cells.Item(17, 9).Value = 1.023843013862435
$ws.Cells.Item(17, 10).Value = 1.023236681640314
$ws.Cells.Item(17, 11).Value = 1.022619490252654
$ws.Cells.Item(17, 12).Value = 1.030078904773804
$ws.Cells.Item(17, 13).Value = 1.037026659793874
$ws.Cells.Item(17, 14).Value = 1.011704350266305
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.017178167073867
$ws.Cells.Item(18, 4).Value = 1.019414859951583
$ws.Cells.Item(18, 5).Value = 1.026914142668853
$ws.Cells.Item(18, 6).Value = 1.033925434066817
$ws.Cells.Item(18, 9).Value = 1.023826118601784
$ws.Cells.Item(18, 10).Value = 1.023408195856578
$ws.Cells.Item(18, 11).Value = 1.022796852535647
$ws.Cells.Item(18, 12).Value = 1.030269627252778
$ws.Cells.Item(18, 13).Value = 1.037256518498098
$ws.Cells.Item(18, 14).Value = 1.011763474262471
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.01726478815901
$ws.Cells.Item(19, 4).Value = 1.019490178028921
$ws.Cells.Item(19, 5).Value = 1.026993907584169
$ws.Cells.Item(19, 6).Value = 1.03401848067386
$ws.Cells.Item(19, 9).Value = 1.023820303822167
$ws.Cells.Item(19, 10).Value = 1.023466678575722
$ws.Cells.Item(19, 11).Value = 1.02285733192662
$ws.Cells.Item(19, 12).Value = 1.03033466361896
$ws.Cells.Item(19, 13).Value = 1.037334906630303
$ws.Cells.Item(19, 14).Value = 1.011783629786303
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.016877440273004
$ws.Cells.Item(20, 4).Value = 1.019153399807599
$ws.Cells.Item(20, 5).Value = 1.026637247034679
$ws.Cells.Item(20, 6).Value = 1.033602447089504
$ws.Cells.Item(20, 9).Value = 1.023846096030192
$ws.Cells.Item(20, 10).Value = 1.023205133269756
$ws.Cells.Item(20, 11).Value = 1.022586867424131
$ws.Cells.Item(20, 12).Value = 1.030043825145954
$ws.Cells.Item(20, 13).Value = 1.036984384826867
$ws.Cells.Item(20, 14).Value = 1.011693472844303
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.015618867949066
$ws.Cells.Item(21, 4).Value = 1.018059588316606
$ws.Cells.Item(21, 5).Value = 1.025478882902058
$ws.Cells.Item(21, 6).Value = 1.032251507085556
$ws.Cells.Item(21, 9).Value = 1.02392621804881
$ws.Cells.Item(21, 10).Value = 1.022354878825014
$ws.Cells.Item(21, 11).Value = 1.021707798543563
$ws.Cells.Item(21, 12).Value = 1.029098621600024
$ws.Cells.Item(21, 13).Value = 1.03584565170029
$ws.Cells.Item(21, 14).Value = 1.011400071138333
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.014828110722074
$ws.Cells.Item(22, 4).Value = 1.017372698487875
$ws.Cells.Item(22, 5).Value = 1.024751471951725
$ws.Cells.Item(22, 6).Value = 1.031403360250182
$ws.Cells.Item(22, 9).Value = 1.023973741382756
$ws.Cells.Item(22, 10).Value = 1.02182033074294
$ws.Cells.Item(22, 11).Value = 1.021155274152923
$ws.Cells.Item(22, 12).Value = 1.028504594244827
$ws.Cells.Item(22, 13).Value = 1.035130332009972
$ws.Cells.Item(22, 14).Value = 1.011215373908602
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.015247253587172
$ws.Cells.Item(23, 4).Value = 1.017736752859941
$ws.Cells.Item(23, 5).Value = 1.025137001161882
$ws.Cells.Item(23, 6).Value = 1.031852861336365
$ws.Cells.Item(23, 9).Value = 1.023948818288725
$ws.Cells.Item(23, 10).Value = 1.022103701023303
$ws.Cells.Item(23, 11).Value = 1.02144816067351
$ws.Cells.Item(23, 12).Value = 1.028819474860521
$ws.Cells.Item(23, 13).Value = 1.035509475331172
$ws.Cells.Item(23, 14).Value = 1.011313306599714
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.016898549919705
$ws.Cells.Item(24, 4).Value = 1.019171751827567
$ws.Cells.Item(24, 5).Value = 1.026656682413228
$ws.Cells.Item(24, 6).Value = 1.033625116891363
$ws.Cells.Item(24, 9).Value = 1.023844704313072
$ws.Cells.Item(24, 10).Value = 1.02321938860577
$ws.Cells.Item(24, 11).Value = 1.02260160821347
$ws.Cells.Item(24, 12).Value = 1.030059676029977
$ws.Cells.Item(24, 13).Value = 1.037003486841567
$ws.Cells.Item(24, 14).Value = 1.011698387960088
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.018817244504749
$ws.Cells.Item(25, 4).Value = 1.020840617249542
$ws.Cells.Item(25, 5).Value = 1.028424108096125
$ws.Cells.Item(25, 6).Value = 1.035687152496668
$ws.Cells.Item(25, 9).Value = 1.02371141742819
$ws.Cells.Item(25, 10).Value = 1.02451426692705
$ws.Cells.Item(25, 11).Value = 1.02394091389106
$ws.Cells.Item(25, 12).Value = 1.031499999575048
$ws.Cells.Item(25, 13).Value = 1.038740036470138
$ws.Cells.Item(25, 14).Value = 1.012144275226148
